# Generate Report for Handoff
# Updates the "Latest Handoff" timestamps for the file that was (re-)handed off:
# 658710e4-ab9e-4a2a-b4b5-4f275914acaf (row 5 on every sheet).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G5").Value = "2016-09-06 15:27:35"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H5").Value = "2016-09-06 15:27:30"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H5").Value = "2016-09-06 15:27:35"
